$wb = $excel.ActiveWorkbook

# --- AddOpportunity sheet: append a 4th data row (Team, Inc / QualSpec Group /
#     Private Equity / Family Office), cloned from row 3 so number formats
#     (the "1000"/"9999" text-formatted columns) and styles carry over, then
#     overwrite the cells that actually differ. ---
$wsOpp = $wb.Worksheets.Item("AddOpportunity")

$wsOpp.Range("A3:AD3").Copy()
$wsOpp.Range("A4:AD4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsOpp.Range("A3:AD3").Copy()
$wsOpp.Range("A4:AD4").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = 0

$wsOpp.Range("A4").Value = "Team, Inc"
$wsOpp.Range("B4").Value = "QualSpec Group"
$wsOpp.Range("R4").Value = "Private Equity"
$wsOpp.Range("S4").Value = "Family Office"
$wsOpp.Columns("S").AutoFit() | Out-Null

# --- CompanyUpdates sheet: append a matching 4th row (same Yes / Source -
#     Engagement / No values as the existing rows). ---
$wsUpd = $wb.Worksheets.Item("CompanyUpdates")

$wsUpd.Range("A3:C3").Copy()
$wsUpd.Range("A4:C4").PasteSpecial(-4104) | Out-Null    # xlPasteAll
$excel.CutCopyMode = 0

# --- Selection / active sheet bookkeeping to match the saved view state ---
$wsUpd.Range("B13").Select() | Out-Null
$wsOpp.Activate() | Out-Null
$wsOpp.Range("V12").Select() | Out-Null
